$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped from 45186 to 45188
# for every data row (rows 2 through 396). Update the value in place so the
# existing cell formatting/style is preserved.
$ws.Range("C2:C396").Value = 45188
